$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (D,E) for the new quarters; existing D:K data shifts to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy cell formatting (number format, font, alignment, etc.) from the shifted F:G
# columns into the newly inserted, blank D:E columns so every row keeps its original
# look (date style for the header rows, "#,##0" style for the data rows).
$ws.Range("F7:G102").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate D:M with the restated / extended financial data for every affected row.
$data = @{}
$data[7] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$data[8] = @(4036000, 4185600, 3854500, 3102800, 3495800, 3486000, 3097900, 2456400, 2648600, 2654500)
$data[9] = @(2316800, 2110100, 1782700, 1469700, 1698700, 1751400, 1570500, 1342800, 1413800, 1346200)
$data[10] = @(1719200, 2075500, 2071800, 1633100, 1797200, 1734600, 1527400, 1113600, 1234800, 1308400)
$data[11] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[12] = @(672100, 581200, 598500, 488900, 549600, 481100, 467200, 412300, 432100, 380100)
$data[13] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[14] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[15] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[16] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[17] = @(3871000, 3528400, 3049800, 2424900, 2787100, 2788500, 2473100, 2164700, 2330800, 2249200)
$data[18] = @(165000, 657200, 804700, 677900, 708700, 697500, 624800, 291700, 317800, 405300)
$data[19] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[20] = @(76300, 1409900, 179900, 363900, 106000, 727500, 163000, 71700, 382700, 243000)
$data[21] = @("NA", "NA", "NA", "NA", 1868400, 1469900, 789600, 507400, "NA", "NA")
$data[22] = @(85200, 71100, 70900, 52200, 59500, 82700, 49600, 47000, 42800, 46500)
$data[23] = @(156100, 1996000, 913600, 989600, 755200, 1342400, 738200, 316400, 657700, 601800)
$data[24] = @(71800, 304700, 161200, 166200, 137900, 162800, 83700, 58900, 58300, 152000)
$data[25] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[26] = @(84300, 1691300, 752400, 823400, 617300, 1179600, 654500, 257500, 599400, 449800)
$data[27] = @(289500, 1839700, 950100, 993500, 620000, 1179700, 655200, 258400, 519400, 451200)
$data[28] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[29] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[30] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[31] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[32] = @(-76300, -1409900, -179900, -363900, -106000, -727500, -163000, -71700, -382700, -243000)
$data[33] = @(289500, 1839700, 950100, 993500, 620000, 1179700, 655200, 258400, 519400, 451200)
$data[34] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[35] = @(289500, 1839700, 950100, 993500, 620000, 1179700, 655200, 258400, 519400, 451200)
$data[38] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$data[39] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[40] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[41] = @(4101800, 2341200, 4283000, 3120300, 1645000, 2279100, 2298000, 2118300, 1585000, 1669400)
$data[42] = @(16566400, 13166300, 14651800, 13272800, 15988400, 16952300, 15088300, 13405600, 11480700, 9726700)
$data[43] = @(1009200, 4270600, 838200, 769500, 4256000, 2937200, 2097300, 1317100, 909700, 659100)
$data[44] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[45] = @(1340100, 996600, 6640200, 7589800, 545700, 845200, 712700, 528100, 532600, 837800)
$data[46] = @(23017500, 20774700, 26413100, 24752400, 22435000, 23013800, 20196200, 17369200, 14507900, 12893100)
$data[47] = @(12577900, 12256100, 9385200, 8837700, 8868800, 7409100, 7222500, 7112800, 7040300, 5997400)
$data[48] = @(2657000, 2403800, 2173800, 2028200, 1851400, 1789500, 1727600, 1620800, 1642500, 1626200)
$data[49] = @(4113500, 3826500, 3366100, 3235900, 3157100, 3072800, 3079100, 3034700, 2794300, 2829000)
$data[50] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[51] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[52] = @(1795900, 1601500, 1451400, 1137400, 1046600, 846100, 745300, 695600, 482800, 643400)
$data[53] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[54] = @(44161800, 40862600, 42789600, 39991600, 37359000, 36131300, 32970800, 29833200, 26467900, 23989200)
$data[55] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[56] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[57] = @(5250900, 5007800, 4548500, 4040000, 4084700, 3708000, 3318400, 3220400, 4167200, 3046900)
$data[58] = @(1484200, 1193100, 2057400, 990900, 6862500, 8547500, 7376600, 5542200, 1424500, 573100)
$data[59] = @(1702400, 1607400, 5445300, 7347300, 1230900, 1238300, 1154800, 1051400, 1113000, 1068600)
$data[60] = @(8437600, 7808300, 12051200, 12378300, 12178100, 13493900, 11849800, 9814000, 6704600, 4688600)
$data[61] = @(8148200, 6116100, 5990700, 6514600, 5314900, 5433100, 5042400, 4981800, 5013100, 5201300)
$data[62] = @(1492700, 1286400, 751500, 517800, 517500, 537200, 551500, 537600, 535500, 521300)
$data[63] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[64] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[65] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[66] = @(19986200, 17163000, 20664600, 21514000, 20240500, 19953600, 18328500, 16146100, 13048500, 11153700)
$data[67] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[68] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[69] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[70] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[71] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[72] = @(19181400, 18877900, 17529200, 16388200, 15186500, 14541100, 13292200, 12592400, 12468300, 11891100)
$data[73] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[74] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[75] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[76] = @(24175500, 23699600, 22125000, 18477600, 17118500, 16177700, 14642300, 13687100, 13419300, 12835500)
$data[77] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[80] = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
$data[81] = @(289500, 1839700, 950100, 993500, 620000, 1179700, 655200, 258400, 519400, 451200)
$data[82] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[83] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[84] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[85] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[86] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[87] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[88] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[89] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[90] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[91] = @("NA", "NA", "NA", "NA", -176500, -192500, -167000, -169800, "NA", "NA")
$data[92] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[93] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[94] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[95] = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
$data[96] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[97] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[98] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[99] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[100] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[101] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$data[102] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item($r, 4 + $i).Value = $v
        }
    }
}

$wb.Save()